$d = $word.ActiveDocument

$replacements = @(
    @{old="44×12=528"; new="90×47=4230"},
    @{old="37×38=1406"; new="83×85=7055"},
    @{old="95×87=8265"; new="14×28=392"},
    @{old="61×56=3416"; new="25×64=1600"},
    @{old="90×67=6030"; new="19×85=1615"},
    @{old="98×68=6664"; new="19×79=1501"},
    @{old="63×42=2646"; new="66×95=6270"},
    @{old="40×41=1640"; new="24×69=1656"},
    @{old="67×23=1541"; new="47×81=3807"},
    @{old="64×50=3200"; new="99×52=5148"},
    @{old="81×50=4050"; new="20×36=720"},
    @{old="27×35=945"; new="28×81=2268"},
    @{old="43×49=2107"; new="49×75=3675"},
    @{old="71×99=7029"; new="88×92=8096"},
    @{old="62×52=3224"; new="12×80=960"},
    @{old="33×78=2574"; new="38×31=1178"},
    @{old="75×59=4425"; new="93×24=2232"},
    @{old="64×18=1152"; new="56×85=4760"},
    @{old="93×30=2790"; new="77×67=5159"},
    @{old="11×47=517"; new="65×78=5070"},
    @{old="50×70=3500"; new="64×87=5568"},
    @{old="50×37=1850"; new="68×91=6188"},
    @{old="39×42=1638"; new="95×49=4655"},
    @{old="28×41=1148"; new="65×14=910"},
    @{old="69×40=2760"; new="58×72=4176"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
